$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 465, shifting rows 465:520 down to 466:521
$ws.Rows("465:465").Insert()

# Populate the newly inserted row 465 with the new data record
$ws.Cells.Item(465, 1).Value = 4
$ws.Cells.Item(465, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(465, 3).Value = "Los Lagos"
$ws.Cells.Item(465, 4).Value = 45194
$ws.Cells.Item(465, 5).Value = 10
$ws.Cells.Item(465, 6).Value = 100112037
$ws.Cells.Item(465, 7).Value = "Cebollín"
$ws.Cells.Item(465, 8).Value = "Sin especificar"
$ws.Cells.Item(465, 9).Value = "Segunda"
$ws.Cells.Item(465, 10).Value = 70
$ws.Cells.Item(465, 11).Value = 6500
$ws.Cells.Item(465, 12).Value = 6500
$ws.Cells.Item(465, 13).Value = 6500
$ws.Cells.Item(465, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(465, 15).Value = "Región Metropolitana"
$ws.Cells.Item(465, 16).Value = 181
$ws.Cells.Item(465, 17).Value = 36
$ws.Cells.Item(465, 18).Value = "Hortaliza"
